$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text updates for columns B (Coin) and C (Link): plain strings, safe to
# assign directly since Excel will not reinterpret them as numbers.
$textUpdates = @{
    'B7' = 'GateToken'
    'C7' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'B9' = 'FTXToken'
    'C9' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'B15' = 'One'
    'C15' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
}

# Numeric-looking updates for columns D (Price) and E (Volume). These must be
# written with a leading apostrophe so Excel keeps them stored as text values
# (the workbook stores them as inline/shared strings, e.g. "255.93", "4.23%"),
# instead of silently converting them into floating point numbers.
$numericTextUpdates = @{
    'D2' = '255.93'
    'E2' = '4.23%'
    'D3' = '28.32'
    'E3' = '-3.51%'
    'D4' = '5.317'
    'E4' = '3.07%'
    'D5' = '0.05798'
    'E5' = '0.53%'
    'D6' = '6.696'
    'E6' = '1.46%'
    'D7' = '3.236'
    'E7' = '2.79%'
    'D8' = '0.8716'
    'E8' = '1.53%'
    'D9' = '0.9128'
    'E9' = '6.02%'
    'D10' = '0.1409'
    'E10' = '3.38%'
    'D11' = '0.07158'
    'E11' = '1.91%'
    'D12' = '0.03160'
    'E12' = '4.40%'
    'D13' = '0.09227'
    'E13' = '-1.45%'
    'D14' = '0.001541'
    'E14' = '-0.17%'
    'D15' = '0.0006081'
    'E15' = '0.88%'
    'D16' = '0.005942'
    'E16' = '-1.77%'
    'D17' = '3.507'
    'E17' = '0.41%'
    'D18' = '2.271'
    'E18' = '5.00%'
    'D19' = '0.3129'
    'E19' = '-2.27%'
    'D20' = '0.03410'
    'E20' = '3.33%'
    'D21' = '0.1312'
    'E21' = '2.37%'
    'D22' = '3.513'
    'E22' = '11.03%'
    'D23' = '0.04160'
    'E24' = '-1.72%'
    'D25' = '0.001221'
    'E25' = '-0.44%'
    'D26' = '0.004977'
    'E26' = '20.32%'
    'D27' = '0.0001198'
    'E27' = '-0.98%'
    'D28' = '0.0001935'
    'E28' = '33.59%'
    'D40' = '0.03874'
    'E40' = '3.80%'
    'D41' = '0.005713'
    'E41' = '-2.98%'
    'E42' = '2.35%'
    'E43' = '-4.80%'
    'D44' = '0.01094'
    'E44' = '30.31%'
    'D45' = '0.00005270'
    'E45' = '-0.20%'
    'D46' = '0.00000000749'
    'E46' = '-0.13%'
    'D47' = '0.08486'
    'E47' = '46.37%'
    'D48' = '0.002169'
    'E48' = '-11.29%'
    'D49' = '0.00002097'
    'E49' = '-0.13%'
    'D50' = '0.0001997'
    'E50' = '-0.13%'
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}

foreach ($cell in $numericTextUpdates.Keys) {
    $ws.Range($cell).Value = "'" + $numericTextUpdates[$cell]
}
